# Apply updated "asked_total" (column D) and "universe" (column J) counts
# to the worker variable map, per the commit "added downloads data
# preperation script" (download totals increased from 210 to 258).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column J ("universe") is 258 for every data row (2-38).
$ws.Range("J2:J38").Value = 258

# Column D ("asked_total") is 258 for every data row (2-38) EXCEPT the
# rows below, whose totals scale differently and row 26 which keeps its
# original value.
$ws.Range("D2:D38").Value = 258

$ws.Range("D3").Value = 23
$ws.Range("D15").Value = 48
$ws.Range("D26").Value = 12
$ws.Range("D32").Value = 132
